$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BUG REPORT")
$ws.Columns.Item(2).Delete()
